$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Update the CasesTab query in B2: append an ORDER BY / LIMIT clause ---
$casesQuery = $ws.Range("B2").Value2
$casesQuery = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B2").Value2 = $casesQuery

# --- Update the SamplesTab query in B3: append an ORDER BY / LIMIT clause ---
$samplesQuery = $ws.Range("B3").Value2
$samplesQuery = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value2 = $samplesQuery

# --- Update the FilesTab query in B4: replace the trailing "order by" clause ---
$filesQuery = $ws.Range("B4").Value2
$filesQuery = $filesQuery.Replace("    order by f.file_name", "  order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value2 = $filesQuery

# --- Row heights grow by one wrapped line now that the queries are longer ---
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# --- Move the active selection to B3 (was B4) ---
$ws.Range("B3").Select()

$wb.Save()
